$wb = $excel.ActiveWorkbook

# Productdata sheet: G2 (AverageDemand) 40 -> 70
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("G2").Value = 70

# Re-blank the empty LostSale (H) column cells so the COM round-trip does
# not materialize them with a stray shared-string value of 0 ("Name") on save.
$wsProductdata.Range("H2:H11").Value = ""

# ForecastedAverageDemand sheet: B9, B10, B11 -> 100
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

# ForcastedStandardDeviation sheet: B9, B10, B11 -> new std dev values
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
